$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing data right by one column
$ws.Columns("A").Insert()

# New column header for the inserted "Match ID" column
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Fill Match ID values for the visible data rows (previously rows 4-19)
$ws.Range("A4:A19").Value = 11
$ws.Range("A4:A19").Font.Bold = $true

# Summary row (row 20) also gets the Match ID value, but without bold formatting
$ws.Range("A20").Value = 11
$ws.Rows(20).EntireRow.AutoFit()

# Restore the selection to the newly added Match ID column data range
$ws.Range("A3:A19").Select() | Out-Null
